$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28, shifting existing rows 28..122 down to 29..123
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new record
$ws.Cells.Item(28,1).Value  = 5
$ws.Cells.Item(28,2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(28,3).Value  = "Maule"
$ws.Cells.Item(28,4).Value  = 44623
$ws.Cells.Item(28,5).Value  = 7
$ws.Cells.Item(28,6).Value  = 100112030
$ws.Cells.Item(28,7).Value  = "Poroto granado"
$ws.Cells.Item(28,8).Value  = "Sin especificar"
$ws.Cells.Item(28,9).Value  = "Primera"
$ws.Cells.Item(28,10).Value = 220
$ws.Cells.Item(28,11).Value = 22000
$ws.Cells.Item(28,12).Value = 23000
$ws.Cells.Item(28,13).Value = 22545
$ws.Cells.Item(28,14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(28,15).Value = "Región del Maule"
$ws.Cells.Item(28,16).Value = 902
$ws.Cells.Item(28,17).Value = 25
$ws.Cells.Item(28,18).Value = "Hortaliza"
